$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at E, shifting the existing Date/Grade columns (and their
# data) right to F/G, to make room for a new "Unnamed: 0.1.1.1" index column.
$ws.Columns("E:E").Insert()

# Header for the newly inserted column
$ws.Range("E1").Value = "Unnamed: 0.1.1.1"

# Fill in the newly revealed index values for existing rows
$ws.Range("D2").Value = 0
$ws.Range("C3").Value = 1
$ws.Range("B4").Value = 2

# Add a new row 5 for the latest game's score
$ws.Range("A4").Copy($ws.Range("A5"))
$ws.Range("A5").Value = 3
$ws.Range("F5").Value = "Sat Jan 18 19:17:08 2020"
$ws.Range("G5").Value = 20
